$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> (new D value, new E value); $null means "leave unchanged"
$updates = @{
    2  = @{ D = "261.03";     E = "0.08%" }
    3  = @{ D = "27.03";      E = "-0.61%" }
    4  = @{ D = "4.710";      E = "-0.39%" }
    5  = @{ D = $null;        E = "2.28%" }
    6  = @{ D = $null;        E = "1.29%" }
    7  = @{ D = "0.8519";     E = "0.50%" }
    8  = @{ D = "0.9157";     E = "-0.52%" }
    9  = @{ D = "0.1404";     E = "-0.29%" }
    10 = @{ D = "0.05060";    E = "2.83%" }
    11 = @{ D = "0.07071";    E = "-0.19%" }
    12 = @{ D = "0.03097";    E = "-1.20%" }
    13 = @{ D = "0.09056";    E = "-0.22%" }
    14 = @{ D = "0.001530";   E = "-1.10%" }
    15 = @{ D = "0.0006159";  E = "0.85%" }
    16 = @{ D = "0.005986";   E = "-3.22%" }
    17 = @{ D = "3.446";      E = "-0.14%" }
    18 = @{ D = "3.169";      E = "0.49%" }
    19 = @{ D = $null;        E = "-1.41%" }
    21 = @{ D = "0.1311";     E = "1.06%" }
    22 = @{ D = "4.101";      E = "0.04%" }
    23 = @{ D = "0.04251";    E = "0.34%" }
    24 = @{ D = "0.001196";   E = "-1.96%" }
    25 = @{ D = "0.004079";   E = "4.27%" }
    27 = @{ D = $null;        E = "4.12%" }
    40 = @{ D = "0.03948";    E = "1.89%" }
    41 = @{ D = "0.1112";     E = "-0.16%" }
    42 = @{ D = $null;        E = "0.17%" }
    43 = @{ D = $null;        E = "0.13%" }
    44 = @{ D = "0.01345";    E = "-17.73%" }
    45 = @{ D = "0.00005162"; E = "-3.19%" }
    46 = @{ D = $null;        E = "0.03%" }
    48 = @{ D = "0.2518";     E = "90.59%" }
    49 = @{ D = $null;        E = "0.03%" }
    50 = @{ D = $null;        E = "0.03%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals.D
        $ws.Range("D$row").Style = "Normal"
    }
    if ($null -ne $vals.E) {
        $ws.Range("E$row").NumberFormat = "@"
        $ws.Range("E$row").Value = $vals.E
        $ws.Range("E$row").Style = "Normal"
    }
}
